# Insert a new price record as row 78, pushing the existing row 78 (and all
# rows below it) down by one. The sheet's dimension grows from A1:R144 to
# A1:R145 automatically as a result of the insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("78:78").Insert()

$ws.Range("A78").Value = 4
$ws.Range("B78").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C78").Value = "Los Lagos"
$ws.Range("D78").Value = "2023-07-11"
$ws.Range("E78").Value = 10
$ws.Range("F78").Value = 100112026
$ws.Range("G78").Value = "Haba"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 80
$ws.Range("K78").Value = 21000
$ws.Range("L78").Value = 21000
$ws.Range("M78").Value = 21000
$ws.Range("N78").Value = "$/saco 25 kilos"
$ws.Range("O78").Value = "Provincia de Limarí"
$ws.Range("P78").Value = 840
$ws.Range("Q78").Value = 25
$ws.Range("R78").Value = "Hortaliza"
